# Regenerate merged AHB files
# - rename the diff-table header labels from *_old/*_new to *_FV2210/*_FV2304
# - turn the sheet's used range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row (row 1) labels.
$headers = @(
    "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210", "Segment ID_FV2210",
    "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210", "Bedingungsausdruck_FV2210", "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
    "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Convert the whole used range (A1:U57) into a native Excel Table, with
#    headers taken from row 1 (xlSrcRange = 1, XlYesNoGuess.xlYes = 1).
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), [Type]::Missing, 1)
$tbl.Name = "Table1"

# 3. Freeze panes so the header row stays put while scrolling.
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true

Write-Output "Regenerated merged AHB sheet: renamed headers, added Table1, froze header row."
